$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Tenant-name labels in the data rows (rows 36-38) ---
$ws.Range("B36").Value = "Мохамед"
$ws.Range("B37").Value = "Банницин"
$ws.Range("B38").Value = "Куропаткин 1"

# --- Summary labels at the bottom of the sheet ---
$ws.Range("B40").Value = "Общая сумма, руб."
$ws.Range("B43").Value = "Максимальная сумма, руб."

# --- View state: scroll the window down so row 25 is at the top and
#     leave the selection on B43 (matches the saved workbook view) ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B43").Select()
